$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.920.82"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").Value = "2.203.61"
$ws.Range("E3").Value = "  -1.74%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.616"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.43%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.41"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.72%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.400"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0893"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.77%  "
$ws.Range("E11").Value = "  -0.43%  "
$ws.Range("D12").Value = "2.532.62"
$ws.Range("E12").Value = "  -1.90%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.34"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.94"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.792"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.99%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.54"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.68%  "
$ws.Range("D17").Value = "2.190.84"
$ws.Range("E17").Value = "  -2.12%  "
$ws.Range("D18").Value = "41.827.61"
$ws.Range("E18").Value = "  +0.30%  "
$ws.Range("D19").Value = "0.0₃0930"
$ws.Range("E19").Value = "  +4.25%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.46%  "
$ws.Range("B21").Value = "Litecoin"
$ws.Range("C21").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "241.75"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.46%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("E25").Value = "  +4.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.56"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "168.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.42%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.139"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.11"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.44%  "
$ws.Range("E30").Value = "  -0.54%  "
$ws.Range("E31").Value = "  -3.81%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.120"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.91"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.57"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.79%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0642"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.81%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.26"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.80%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.51"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.85%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.31"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0245"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.77%  "
$ws.Range("E40").Value = "  -0.18%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.46"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.000219"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -10.23%  "
$ws.Range("E43").Value = "  -2.87%  "
$ws.Range("E44").Value = "  +0.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "96.25"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.39%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "1.451.83"
$ws.Range("E46").Value = "  -1.81%  "
$ws.Range("B47").Value = "FTXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.31"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -13.30%  "
$ws.Range("E48").Value = "  -1.56%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "16.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.45%  "
$ws.Range("E50").Value = "  -2.90%  "
$ws.Range("E51").Value = "  +0.51%  "
